$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 124.166664
$ws.Range("I4").Value = 124.166664
$ws.Range("K4").Value = 124.166664
$ws.Range("M4").Value = -10.166664

$ws.Range("H74").Value = 73664.14
$ws.Range("I74").Value = 2608.1667
$ws.Range("J74").Value = 500000
$ws.Range("K74").Value = 2608.1667
$ws.Range("L74").Value = 500000
$ws.Range("M74").Value = -1672.1667
$ws.Range("N74").Value = -501872

$ws.Range("H77").Value = 73664.14
$ws.Range("I77").Value = 2608.1667
$ws.Range("J77").Value = 500000
$ws.Range("K77").Value = 13040.8335
$ws.Range("L77").Value = 2500000
$ws.Range("M77").Value = -8360.833500000001
$ws.Range("N77").Value = -2509360

$ws.Range("H104").Value = 205.66667
$ws.Range("I104").Value = 205.66667
$ws.Range("K104").Value = 617.00001
$ws.Range("M104").Value = 1129.99999

$ws.Range("H107").Value = 798.28125
$ws.Range("I107").Value = 721.1923
$ws.Range("K107").Value = 721.1923
$ws.Range("M107").Value = 1198.8077

$ws.Range("H113").Value = 5845.5
$ws.Range("I113").Value = 4593.5713
$ws.Range("K113").Value = 4593.5713
$ws.Range("M113").Value = -1339.5713

$ws.Range("H135").Value = 1945.25
$ws.Range("I135").Value = 1882
$ws.Range("J135").Value = 1966.3334
$ws.Range("K135").Value = 16938
$ws.Range("L135").Value = 17697.0006
$ws.Range("M135").Value = -14403
$ws.Range("N135").Value = -22767.0006

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2444.25
$ws.Range("I2").Value = 1735.909
$ws.Range("K2").Value = 1735.909
$ws.Range("M2").Value = -1622.909

$ws.Range("H116").Value = 2444.25
$ws.Range("I116").Value = 1735.909
$ws.Range("K116").Value = 1735.909
$ws.Range("M116").Value = 558.0909999999999

$ws.Range("H132").Value = 1206.2858
$ws.Range("I132").Value = 1206.2858
$ws.Range("K132").Value = 3618.8574
$ws.Range("M132").Value = -1088.8574

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2444.25
$ws.Range("I3").Value = 1735.909
$ws.Range("K3").Value = 1735.909
$ws.Range("M3").Value = -1621.909

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws.Range("H94").Value = 1477.3914
$ws.Range("I94").Value = 1098.5
$ws.Range("K94").Value = 1098.5
$ws.Range("M94").Value = -647.5

$ws.Range("H95").Value = 18569.4
$ws.Range("J95").Value = 18569.4
$ws.Range("L95").Value = 18569.4
$ws.Range("N95").Value = -24061.4

$ws.Range("H107").Value = 3173.5
$ws.Range("I107").Value = 1347
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 1347
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = 573
$ws.Range("N107").Value = -8840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4637.1177
$ws.Range("I31").Value = 4276.8
$ws.Range("K31").Value = 4276.8
$ws.Range("M31").Value = -3981.8

$ws.Range("H34").Value = 4637.1177
$ws.Range("I34").Value = 4276.8
$ws.Range("K34").Value = 4276.8
$ws.Range("M34").Value = -4074.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H34").Value = 1083.8334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 31000
$ws.Range("J49").Value = 31000
$ws.Range("L49").Value = 31000
$ws.Range("N49").Value = -31368

$ws.Range("H80").Value = 3722.7778
$ws.Range("I80").Value = 3299.6667
$ws.Range("J80").Value = 3934.3333
$ws.Range("K80").Value = 3299.6667
$ws.Range("L80").Value = 3934.3333
$ws.Range("M80").Value = -2301.6667
$ws.Range("N80").Value = -5930.3333

$ws.Range("H83").Value = 3722.7778
$ws.Range("I83").Value = 3299.6667
$ws.Range("J83").Value = 3934.3333
$ws.Range("K83").Value = 16498.3335
$ws.Range("L83").Value = 19671.6665
$ws.Range("M83").Value = -11506.3335
$ws.Range("N83").Value = -29655.6665

$ws.Range("H97").Value = 918.58826
$ws.Range("I97").Value = 932.0769
$ws.Range("K97").Value = 932.0769
$ws.Range("M97").Value = -436.0769

$ws.Range("H102").Value = 3107.7144
$ws.Range("I102").Value = 3107.7144
$ws.Range("K102").Value = 3107.7144
$ws.Range("M102").Value = -1485.7144

$ws.Range("H132").Value = 797.3333
$ws.Range("I132").Value = 797.3333
$ws.Range("K132").Value = 2391.9999
$ws.Range("M132").Value = 138.0001000000002

$ws.Range("H133").Value = 87832.5
$ws.Range("J133").Value = 87832.5
$ws.Range("L133").Value = 87832.5
$ws.Range("N133").Value = -97952.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2120.8333
$ws.Range("J22").Value = 2272.7273
$ws.Range("L22").Value = 2272.7273
$ws.Range("N22").Value = -2862.7273

$ws.Range("H27").Value = 2120.8333
$ws.Range("J27").Value = 2272.7273
$ws.Range("L27").Value = 2272.7273
$ws.Range("N27").Value = -2486.7273

$ws.Range("H45").Value = 21000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 21000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 21000
$ws.Range("N45").Value = -21814
$ws.Range("M45").ClearContents()

$ws.Range("H46").Value = 3095.6924
$ws.Range("I46").Value = 2500
$ws.Range("K46").Value = 2500
$ws.Range("M46").Value = -2312

$ws.Range("H68").Value = 4000
$ws.Range("I68").Value = 4000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 4000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -3251
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 4000
$ws.Range("I71").Value = 4000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 20000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -16256
$ws.Range("N71").ClearContents()

$ws.Range("H74").Value = 29000
$ws.Range("I74").Value = 18000
$ws.Range("J74").Value = 40000
$ws.Range("K74").Value = 18000
$ws.Range("L74").Value = 40000
$ws.Range("M74").Value = -17002
$ws.Range("N74").Value = -41996

$ws.Range("H77").Value = 29000
$ws.Range("I77").Value = 18000
$ws.Range("J77").Value = 40000
$ws.Range("K77").Value = 54000
$ws.Range("L77").Value = 120000
$ws.Range("M77").Value = -49008
$ws.Range("N77").Value = -129984

$ws.Range("H93").Value = 1299.6666
$ws.Range("I93").Value = 1299.6666
$ws.Range("K93").Value = 1299.6666
$ws.Range("M93").Value = -51.66660000000002

$ws.Range("H122").Value = 4388
$ws.Range("I122").Value = 4271
$ws.Range("K122").Value = 12813
$ws.Range("M122").Value = -10363

$ws.Range("H132").Value = 9713.385
$ws.Range("I132").Value = 9435.944
$ws.Range("K132").Value = 28307.832
$ws.Range("M132").Value = -25777.832

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 323.25
$ws.Range("I107").Value = 326.57144
$ws.Range("K107").Value = 979.71432
$ws.Range("M107").Value = 940.28568

$ws.Range("H110").Value = 61052.5
$ws.Range("J110").Value = 61052.5
$ws.Range("L110").Value = 61052.5
$ws.Range("N110").Value = -69232.5

$ws.Range("H113").Value = 1981.8
$ws.Range("I113").Value = 977.25
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 2931.75
$ws.Range("L113").Value = 18000
$ws.Range("M113").Value = -761.75
$ws.Range("N113").Value = -22340

$ws.Range("H132").Value = 3114.5
$ws.Range("I132").Value = 2486
$ws.Range("K132").Value = 7458
$ws.Range("M132").Value = -4928

Write-Output "Done applying Marilith_Profits edits"
